$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (avoid numeric auto-coercion),
# matching the original inlineStr string cells, then restore default styling.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.011.01"
$ws.Range("D3").Value = "1.908.47"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "319.68"
$ws.Range("D7").Value = "0.5045"
$ws.Range("D8").Value = "0.4053"
$ws.Range("D9").Value = "0.08282"
$ws.Range("D10").Value = "42.01"
$ws.Range("D11").Value = "1.100"
$ws.Range("D12").Value = "24.27"
$ws.Range("D13").Value = "1.910.70"
$ws.Range("D14").Value = "6.386"
$ws.Range("D15").Value = "7.245"
$ws.Range("D16").Value = "1.001"
$ws.Range("D17").Value = "91.90"
$ws.Range("D18").Value = "0.00001095"
$ws.Range("D19").Value = "0.06508"
$ws.Range("D20").Value = "18.06"
$ws.Range("D21").Value = "1.000"
$ws.Range("D22").Value = "5.935"
$ws.Range("D23").Value = "30.038.39"
$ws.Range("D24").Value = "11.28"
$ws.Range("D25").Value = "2.204"
$ws.Range("D26").Value = "22.25"
$ws.Range("D27").Value = "2.129.51"
$ws.Range("D28").Value = "161.84"
$ws.Range("D29").Value = "2.275"
$ws.Range("D30").Value = "128.82"
$ws.Range("D31").Value = "1.118"
$ws.Range("D33").Value = "5.949"
$ws.Range("D34").Value = "3.806"
$ws.Range("D35").Value = "5.401"
$ws.Range("D36").Value = "0.02436"
$ws.Range("D37").Value = "0.06340"
$ws.Range("D38").Value = "0.2151"
$ws.Range("D39").Value = "0.6571"
$ws.Range("D40").Value = "1.192"
$ws.Range("D41").Value = "8.677"
$ws.Range("D42").Value = "11.35"
$ws.Range("D43").Value = "1.205"
$ws.Range("D44").Value = "2.206"
$ws.Range("D45").Value = "13.26"
$ws.Range("D46").Value = "0.6045"
$ws.Range("D47").Value = "3.636"
$ws.Range("D48").Value = "122.98"
$ws.Range("D50").Value = "78.32"
$ws.Range("D51").Value = "1.129"

$ws.Range("D2:D51").Style = "Normal"

$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  +2.54%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  -5.34%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("E42").Value = "  -4.90%  "
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("E44").Value = "  +6.96%  "
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("E51").Value = "  -2.71%  "

